$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new data rows above the current row 305, shifting existing rows
# 305-364 down to 308-367 (matches the dimension change A1:T364 -> A1:T367).
$ws.Rows("305:307").Insert()

# Row 305 (new)
$ws.Range("A305").Value = 10
$ws.Range("B305").Value = "Vega Modelo de Temuco"
$ws.Range("C305").Value = "La Araucanía"
$ws.Range("D305").Value = 44508
$ws.Range("E305").Value = 9
$ws.Range("F305").Value = "Fruta"
$ws.Range("G305").Value = 100108
$ws.Range("H305").Value = "Tropicales y subtropicales"
$ws.Range("I305").Value = 100108006
$ws.Range("J305").Value = "Plátano"
$ws.Range("K305").Value = "Barraganete"
$ws.Range("L305").Value = "Primera"
$ws.Range("M305").Value = 65
$ws.Range("N305").Value = 30000
$ws.Range("O305").Value = 30000
$ws.Range("P305").Value = 30000
$ws.Range("Q305").Value = "$/caja 20 kilos"
$ws.Range("R305").Value = "Ecuador"
$ws.Range("S305").Value = 1500
$ws.Range("T305").Value = 20

# Row 306 (new)
$ws.Range("A306").Value = 10
$ws.Range("B306").Value = "Vega Modelo de Temuco"
$ws.Range("C306").Value = "La Araucanía"
$ws.Range("D306").Value = 44508
$ws.Range("E306").Value = 9
$ws.Range("F306").Value = "Fruta"
$ws.Range("G306").Value = 100108
$ws.Range("H306").Value = "Tropicales y subtropicales"
$ws.Range("I306").Value = 100108006
$ws.Range("J306").Value = "Plátano"
$ws.Range("K306").Value = "Sin especificar"
$ws.Range("L306").Value = "Maduro"
$ws.Range("M306").Value = 125
$ws.Range("N306").Value = 17000
$ws.Range("O306").Value = 17000
$ws.Range("P306").Value = 17000
$ws.Range("Q306").Value = "$/caja 20 kilos"
$ws.Range("R306").Value = "Ecuador"
$ws.Range("S306").Value = 850
$ws.Range("T306").Value = 20

# Row 307 (new)
$ws.Range("A307").Value = 10
$ws.Range("B307").Value = "Vega Modelo de Temuco"
$ws.Range("C307").Value = "La Araucanía"
$ws.Range("D307").Value = 44508
$ws.Range("E307").Value = 9
$ws.Range("F307").Value = "Fruta"
$ws.Range("G307").Value = 100108
$ws.Range("H307").Value = "Tropicales y subtropicales"
$ws.Range("I307").Value = 100108006
$ws.Range("J307").Value = "Plátano"
$ws.Range("K307").Value = "Sin especificar"
$ws.Range("L307").Value = "Pintón"
$ws.Range("M307").Value = 1530
$ws.Range("N307").Value = 18000
$ws.Range("O307").Value = 22000
$ws.Range("P307").Value = 20301
$ws.Range("Q307").Value = "$/caja 20 kilos"
$ws.Range("R307").Value = "Ecuador"
$ws.Range("S307").Value = 1015
$ws.Range("T307").Value = 20
